$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.283.94'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").Value = '1.866.47'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.08%  '
$c = $ws.Range("D5")
$c.Value = "'234.83"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("E6").Value = '  -0.04%  '
$c = $ws.Range("D7")
$c.Value = "'0.4697"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -0.30%  '
$c = $ws.Range("D8")
$c.Value = "'0.2857"
$c.Style = "Normal"
$ws.Range("E8").Value = '  -1.69%  '
$c = $ws.Range("D9")
$c.Value = "'0.06567"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '
$c = $ws.Range("D10")
$c.Value = "'21.35"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.25%  '
$c = $ws.Range("D11")
$c.Value = "'0.07829"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.38%  '
$c = $ws.Range("D12")
$c.Value = "'96.75"
$c.Style = "Normal"
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("D13").Value = '1.854.81'
$ws.Range("E13").Value = '  -0.42%  '
$c = $ws.Range("D14")
$c.Value = "'0.6969"
$c.Style = "Normal"
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("E15").Value = '  -1.07%  '
$c = $ws.Range("D16")
$c.Value = "'268.13"
$c.Style = "Normal"
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").Value = '30.405.32'
$ws.Range("E17").Value = '  +0.49%  '
$c = $ws.Range("D18")
$c.Value = "'13.82"
$c.Style = "Normal"
$ws.Range("E18").Value = '  +0.39%  '
$c = $ws.Range("D19")
$c.Value = "'0.000007652"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +2.68%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '2.134.16'
$ws.Range("E21").Value = '  +1.43%  '
$ws.Range("E22").Value = '  -0.03%  '
$c = $ws.Range("D23")
$c.Value = "'5.235"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -0.62%  '
$c = $ws.Range("D24")
$c.Value = "'6.174"
$c.Style = "Normal"
$ws.Range("E24").Value = '  +0.02%  '
$c = $ws.Range("D25")
$c.Value = "'9.458"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +2.94%  '
$c = $ws.Range("D26")
$c.Value = "'166.66"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -0.41%  '
$c = $ws.Range("D27")
$c.Value = "'18.86"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("E28").Value = '  -0.69%  '
$c = $ws.Range("D29")
$c.Value = "'1.366"
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.02%  '
$c = $ws.Range("D30")
$c.Value = "'0.09913"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +0.65%  '
$c = $ws.Range("D31")
$c.Value = "'4.362"
$c.Style = "Normal"
$ws.Range("E31").Value = '  +0.25%  '
$c = $ws.Range("D32")
$c.Value = "'1.457"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.92%  '
$c = $ws.Range("D33")
$c.Value = "'4.045"
$c.Style = "Normal"
$ws.Range("E33").Value = '  +0.46%  '
$c = $ws.Range("D34")
$c.Value = "'0.04719"
$c.Style = "Normal"
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("E35").Value = '  +0.15%  '
$c = $ws.Range("D36")
$c.Value = "'0.7021"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +0.27%  '
$ws.Range("E37").Value = '  +0.30%  '
$c = $ws.Range("D38")
$c.Value = "'0.01875"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '
$c = $ws.Range("D39")
$c.Value = "'2.755"
$c.Style = "Normal"
$ws.Range("E39").Value = '  +5.04%  '
$c = $ws.Range("D40")
$c.Value = "'6.324"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.44%  '
$c = $ws.Range("D41")
$c.Value = "'72.97"
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.37%  '
$c = $ws.Range("D42")
$c.Value = "'1.948"
$c.Style = "Normal"
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("E43").Value = '  +0.42%  '
$ws.Range("E44").Value = '  +0.00%  '
$c = $ws.Range("D45")
$c.Value = "'0.8362"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -0.82%  '
$c = $ws.Range("D46")
$c.Value = "'102.98"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.29%  '
$c = $ws.Range("D47")
$c.Value = "'971.99"
$c.Style = "Normal"
$ws.Range("E47").Value = '  +2.62%  '
$c = $ws.Range("D48")
$c.Value = "'7.108"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.72%  '
$c = $ws.Range("D49")
$c.Value = "'9.116"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -1.01%  '
$c = $ws.Range("D50")
$c.Value = "'34.46"
$c.Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  +0.33%  '
